$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.263.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.05%  "
$ws.Range("D3").Value = "'1.783.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "'338.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.3817"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.46%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "'46.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -2.82%  "
$ws.Range("D11").Value = "'0.07368"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "'23.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.30%  "
$ws.Range("D13").Value = "'0.9993"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "'6.446"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").Value = "'7.386"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.21%  "
$ws.Range("D16").Value = "'1.790.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "'0.00001075"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "'0.06674"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "
$ws.Range("D19").Value = "'82.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'6.453"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'28.258.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.04%  "
$ws.Range("D24").Value = "'12.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").Value = "'2.371"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'1.455"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.15%  "
$ws.Range("D27").Value = "'20.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Value = "'2.415"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "
$ws.Range("D29").Value = "'154.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'136.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "
$ws.Range("D31").Value = "'1.985.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "
$ws.Range("D32").Value = "'6.112"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").Value = "'3.947"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("D34").Value = "'0.08887"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.16%  "
$ws.Range("D35").Value = "'12.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "'0.02441"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("D37").Value = "'0.6855"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "'5.336"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").Value = "'0.06359"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "'0.2173"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "'1.245"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("D42").Value = "'1.494"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.77%  "
$ws.Range("D43").Value = "'8.316"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "'14.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'0.9996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").Value = "'0.6296"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "'3.877"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'133.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("D49").Value = "'2.086"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("D50").Value = "'0.07439"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.61%  "
$ws.Range("D51").Value = "'1.205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.06%  "
